# Pinkbike BuySell All Mountain template update
# The "MySQL for Excel" add-in relabelled the generic HTML-scraping
# terminology used in column A of Sheet1:
#   element_name            -> html_tag
#   element_attribute_name  -> html_tag_attribute_name
#   element_attribute_value -> html_tag_attribute_value
#   Attribute_Name          -> Item_Attribute_Name
# The values in column B (div, class, bsitem, style, ...) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Relabel column A, top to bottom, so new shared-string entries are
# --- created in the same first-seen order as the authoritative edit.
$ws.Range("A6").Value  = "html_tag"
$ws.Range("A7").Value  = "html_tag_attribute_name"
$ws.Range("A8").Value  = "html_tag_attribute_value"

$ws.Range("A11").Value = "Item_Attribute_Name"
$ws.Range("A12").Value = "html_tag"
$ws.Range("A13").Value = "html_tag_attribute_name"
$ws.Range("A14").Value = "html_tag_attribute_value"

$ws.Range("A17").Value = "Item_Attribute_Name"
$ws.Range("A18").Value = "html_tag"
$ws.Range("A19").Value = "html_tag_attribute_name"
$ws.Range("A20").Value = "html_tag_attribute_value"

$ws.Range("A25").Value = "Item_Attribute_Name"
$ws.Range("A26").Value = "html_tag"
$ws.Range("A27").Value = "html_tag_attribute_name"
$ws.Range("A28").Value = "html_tag_attribute_value"

$ws.Range("A30").Value = "Item_Attribute_Name"
$ws.Range("A31").Value = "html_tag"
$ws.Range("A32").Value = "html_tag_attribute_name"
$ws.Range("A33").Value = "html_tag_attribute_value"

$ws.Range("A35").Value = "Item_Attribute_Name"
$ws.Range("A36").Value = "html_tag"
$ws.Range("A37").Value = "html_tag_attribute_name"
$ws.Range("A38").Value = "html_tag_attribute_value"

# --- Workbook-level defined name added by the add-in refresh
# --- (hidden helper formula used to build MySQL-style datetime formats).
$dateFormatFormula = "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)"
$definedName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", $dateFormatFormula)
$definedName.Visible = $false

# --- Leave the user's selection where the edit session ended.
$ws.Range("B30").Select()
